# This workbook contains a weekly price table for "Zapallo" (Camote / Paine
# varieties) reported by "Macroferia Regional de Talca". The commit adds two
# new weekly observations. In the OOXML this shows up as two brand-new rows
# being inserted right after row 23 (pushing every following row down by two,
# from row 24..131 to row 26..133), and the two new rows being populated with
# the new observations' data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows at position 24, pushing old rows 24:131 down to 26:133.
$ws.Rows("24:25").Insert()

# Populate the first new row (new weekly "Camote" observation).
$ws.Range("A24").Value = 5
$ws.Range("B24").Value = "Macroferia Regional de Talca"
$ws.Range("C24").Value = "Maule"
$ws.Range("D24").Value = 44453
$ws.Range("E24").Value = 7
$ws.Range("F24").Value = 100112045
$ws.Range("G24").Value = "Zapallo"
$ws.Range("H24").Value = "Camote"
$ws.Range("I24").Value = "1a (guarda)"
$ws.Range("J24").Value = 900
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 600
$ws.Range("N24").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O24").Value = "Región del Maule"
$ws.Range("P24").Value = 600
$ws.Range("Q24").Value = 1
$ws.Range("R24").Value = "Hortaliza"

# Populate the second new row (new weekly "Paine" observation).
$ws.Range("A25").Value = 5
$ws.Range("B25").Value = "Macroferia Regional de Talca"
$ws.Range("C25").Value = "Maule"
$ws.Range("D25").Value = 44453
$ws.Range("E25").Value = 7
$ws.Range("F25").Value = 100112045
$ws.Range("G25").Value = "Zapallo"
$ws.Range("H25").Value = "Paine"
$ws.Range("I25").Value = "1a (guarda)"
$ws.Range("J25").Value = 2000
$ws.Range("K25").Value = 150
$ws.Range("L25").Value = 150
$ws.Range("M25").Value = 150
$ws.Range("N25").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O25").Value = "Región del Maule"
$ws.Range("P25").Value = 150
$ws.Range("Q25").Value = 1
$ws.Range("R25").Value = "Hortaliza"
